$d = $word.ActiveDocument

# 1. Update the letter date: "September 19, 2025" -> "September 21, 2025"
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.Trim() -eq "September 19, 2025") {
        $para.Range.Text = "September 21, 2025"
        break
    }
}

# 2. Split the mailing-address paragraph "2525 DEVELOPER, Santa Clara CA 95070"
#    into two paragraphs: "2525 DEVELOPER" and "Santa Clara, CA 95070".
#    Only the recipient address block (outside any table) changes - the
#    "PROPERTY ADDRESS:" value inside the info table keeps its original text.
foreach ($para in $d.Paragraphs) {
    $isInTable = $para.Range.Information(12)
    if (-not $isInTable -and $para.Range.Text.Trim() -eq "2525 DEVELOPER, Santa Clara CA 95070") {
        $para.Range.Text = "2525 DEVELOPER"
        $para.Range.InsertParagraphAfter()
        $newPara = $para.Next()
        $newPara.Range.Text = "Santa Clara, CA 95070"
        break
    }
}

# 3. Remove the empty "NoSpacing" paragraph that follows "Board of Directors"
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -match "Board of Directors") {
        $para.Next().Range.Delete()
        break
    }
}
